$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number (e.g. "300.57").
# Force them to Text format first so Excel stores the exact string instead of
# auto-converting to a floating point number, then restore the default style
# so no visible formatting change is introduced.
$riskyCells = @("D5","D6","D7","D9","D10","D11","D12","D16","D21","D23","D24","D28","D29","D31","D32","D33","D34","D37","D38","D39","D45","D46","D47","D49","D50")
foreach ($addr in $riskyCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the crypto price refresh
$ws.Range("D2").Value = "42.629.51"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "2.287.65"
$ws.Range("E3").Value = "  -2.78%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "300.57"
$ws.Range("E5").Value = "  -2.44%  "
$ws.Range("D6").Value = "96.99"
$ws.Range("E6").Value = "  -5.47%  "
$ws.Range("D7").Value = "0.504"
$ws.Range("E7").Value = "  -1.30%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.501"
$ws.Range("E9").Value = "  -3.09%  "
$ws.Range("D10").Value = "33.39"
$ws.Range("E10").Value = "  -5.48%  "
$ws.Range("D11").Value = "0.0789"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "50.14"
$ws.Range("E12").Value = "  -4.52%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("D15").Value = "2.641.93"
$ws.Range("E15").Value = "  -3.06%  "
$ws.Range("D16").Value = "15.43"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "2.299.62"
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("E18").Value = "  -1.89%  "
$ws.Range("D19").Value = "42.555.94"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("E20").Value = "  -1.37%  "
$ws.Range("D21").Value = "11.53"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("E22").Value = "  -4.30%  "
$ws.Range("D23").Value = "66.88"
$ws.Range("E23").Value = "  -1.22%  "
$ws.Range("D24").Value = "234.84"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("E25").Value = "  -3.49%  "
$ws.Range("E26").Value = "  -3.45%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "24.45"
$ws.Range("E28").Value = "  -3.98%  "
$ws.Range("D29").Value = "166.35"
$ws.Range("E29").Value = "  +3.31%  "
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").Value = "33.89"
$ws.Range("E31").Value = "  -5.61%  "
$ws.Range("D32").Value = "9.12"
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").Value = "4.97"
$ws.Range("E34").Value = "  -3.92%  "
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("E36").Value = "  -4.78%  "
$ws.Range("D37").Value = "4.34"
$ws.Range("E37").Value = "  -6.61%  "
$ws.Range("D38").Value = "2.84"
$ws.Range("E38").Value = "  -6.81%  "
$ws.Range("D39").Value = "16.27"
$ws.Range("E39").Value = "  -9.03%  "
$ws.Range("E40").Value = "  -6.60%  "
$ws.Range("E41").Value = "  -3.58%  "
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("E43").Value = "  -3.99%  "
$ws.Range("D44").Value = "1.964.54"
$ws.Range("E44").Value = "  -3.50%  "
$ws.Range("D45").Value = "0.0283"
$ws.Range("E45").Value = "  -1.60%  "
$ws.Range("D46").Value = "17.79"
$ws.Range("E46").Value = "  -7.96%  "
$ws.Range("D47").Value = "9.71"
$ws.Range("E47").Value = "  -7.85%  "
$ws.Range("E48").Value = "  -6.61%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "53.22"
$ws.Range("E49").Value = "  -7.30%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "2.82"
$ws.Range("E50").Value = "  -3.20%  "
$ws.Range("D51").Value = "2.508.37"
$ws.Range("E51").Value = "  -3.22%  "

# Restore default (unstyled) formatting on the cells we touched above
foreach ($addr in $riskyCells) {
    $ws.Range($addr).Style = "Normal"
}
